$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: column B now holds "Телефон" (was "ФИО"); columns C/D headers removed.
$ws.Range("B1").Value = "Телефон"

# Update data row 2: A2 new id, B2 now holds the phone number (was the name).
$ws.Range("A2").Value = 1149518006
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "+73452598001"
$ws.Range("B2").Style = "Normal"

# Remove the now-unused columns C and D entirely (address + old phone columns).
$ws.Range("C1:D2").EntireColumn.Delete()
